$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(1.0344335891558327, 0.6006552548631346, 1.4747487308486456, 3.8848691198726564)"
$ws.Range("C2").Value = "JSU(-1.3562912473480124, 1.4174043453827725, 2.220043201918458, 5.7391831239407605)"
$ws.Range("D2").Value = "JSU(-1.1197051447560884, 1.1352471193249505, -0.06373850907726558, 2.6641709470701467)"
$ws.Range("E2").Value = "JSU(-0.9030508801671311, 1.1654564947743757, 3.560944833938579, 5.512241708039925)"
